# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C85").Value = "[name=`"Suzuran`"]'Dusk Beauty Flying to the Moon'.`n"
$ws.Range("D86").Value = "[name=`"스즈란`"]음… 석아라는 글자는 염국식으로는 '시어'라고 읽는대요. 석아분월은 달로 간 시어라는 뜻이죠. `n"
$ws.Range("C87").Value = "[name=`"Suzuran`"]'Once upon a time, there was a gentle and kind girl named Dusk Beauty who lived in Yen——'`n"
$ws.Range("D87").Value = "[name=`"스즈란`"]옛날 옛적에, 염국에는 시어라는 이름을 가진 마음씨 따뜻한 여자아이가 있었습니다……`n"
$ws.Range("C94").Value = "[name=`"Suzuran`"]'——To commemorate the distant Dusk Beauty, on the same day every year, people would offer up her favorite food. As time passed, this day became an important festival in Yen.'  `n"
$ws.Range("D94").Value = "[name=`"스즈란`"]……사람들은 떠나버린 시어를 기억하기 위해, 매년 이날에 그녀가 가장 좋아하는 음식을 차려놓고 그녀를 기렸습니다. 그렇게 오랜 시간이 지나, 이날은 이곳의 중요한 기념일이 되었습니다……`n"
$ws.Range("D118").Value = "[name=`"스즈란`"]이건 염국의 신화 이야기인데, 주인공은 시어라는 이름을 가진 언니에요.`n"
$ws.Range("D119").Value = "[name=`"스즈란`"]시어는 원래 남편과 함께 행복하게 살고 있었어요.`n"
$ws.Range("D120").Value = "[name=`"스즈란`"]하지만 어느 날, 시어의 남편이 실종되었어요.`n"
$ws.Range("D122").Value = "[name=`"스즈란`"]네, 책에서는 시어가 곳곳을 돌아다니며 남편의 행방을 찾는다고 나오는데, 남편의 행방을 아는 사람은 아무도 없었대요.`n"
$ws.Range("D123").Value = "[name=`"스즈란`"]어떻게 할 방법이 없었던 시어는 그대로 집에 돌아올 수밖에 없었고, 남편이 직접 심은 버드나무 아래에서 두 손으로 얼굴을 가리고 펑펑 울었어요.`n"
$ws.Range("D124").Value = "[name=`"스즈란`"]이때, 지혜로운 마을의 어느 어르신 한 분이, 시어가 펑펑 우는 것을 보고는 좋은 방법을 알려줬어요.`n"
$ws.Range("C127").Value = "[name=`"Suzuran`"]'Child, you have a pair of sharp eyes that can see thousands of miles away. Why not ascend the mountain to the west?' `n"
$ws.Range("C128").Value = "[name=`"Suzuran`"]'From there you can gaze upon the wide earth. Only there will you have a chance to find your husband.' `n"
$ws.Range("D129").Value = "[name=`"스즈란`"]시어는 어르신 말씀이 일리가 있다고 생각해서, 마을 사람들과 작별하고 대지의 끝자락인 고산으로 향했어요.`n"
$ws.Range("C132").Value = "[name=`"Suzuran`"]The meaning of Mount Sky is 'a mountain as tall as the skies'. There is no mountain taller than it in all of Yen. `n"
$ws.Range("D133").Value = "[name=`"스즈란`"]갖은 고난을 이겨내고 산봉우리에 도착한 시어는 온 세상을 널리 볼 수 있었지만, 3일 내내 눈이 새빨개지도록 찾아봐도 남편의 종적은 찾을 수 없었어요.`n"
$ws.Range("D134").Value = "[name=`"스즈란`"]수년간의 노력이 물거품이 되자 시어의 눈물은 쉴 새 없이 쏟아졌고, 눈물은 산맥을 따라 흘러 땅에 떨어져 강물이 되었다고 해요.`n"
$ws.Range("D139").Value = "[name=`"스즈란`"]자초지종을 알게 된 산신들은 더욱 의아해했어요. 그래서 산신들은 시어를 천악산의 신, '천악의 신' 앞에 데려가 이야기를 나누게 했죠.`n"
$ws.Range("D140").Value = "[name=`"스즈란`"]마음씨 좋은 천악의 신은 시어의 사정을 이해하고 가엾게 여겨 시어에게 이렇게 말했어요.`n"
$ws.Range("C141").Value = "[name=`"Suzuran`"]'Lofty heights of Mount Sky, lie yet dwarfed by the twin moons. Ascend thee to the moons, and thy husband shall thee find soon.' `n"
$ws.Range("D146").Value = "[name=`"스즈란`"]시어는 산신의 제안을 받아들여 도움을 청했어요. 머리를 맞대고 고민하던 산신들은 시어에게 방법을 하나 가르쳐 주었죠……`n"
$ws.Range("D153").Value = "[name=`"스즈란`"]그렇지만 어르신은 시어의 일을 떠올리고는, 모두에게 이것은 시어가 남편의 행방을 찾았다는 것이니 무서워할 게 아니라 오히려 축하해야 할 일이라고 말했죠.`n"
$ws.Range("D154").Value = "[name=`"스즈란`"]다들 어르신의 말이 일리가 있다고 생각해서, 집에서 음식과 술을 쉴 새 없이 꺼내와 마을을 떠난 시어를 축하했어요.`n"
$ws.Range("D155").Value = "[name=`"스즈란`"]그래서 해마다 다들 음식을 가지고 한데 모여 먹고 마셨고, 시어가 좋아하는 음식을 상에 올렸대요.`n"
$ws.Range("D158").Value = "[name=`"스즈란`"]시어가 떠난 후에도 강물은 천악산에서 세차게 흘러내렸고, 현지 사람들은 시어의 이름을 따서 이 강물을 석강이라고 불렀대요.`n"
$ws.Range("C184").Value = "[name=`"Suzuran`"]'As long as you can bring me one treasure, I will become your wife.' `n"
$ws.Range("C195").Value = "[name=`"Mousse`"]'As long as you can bring me one treasure, I will become your wife.' `n"
$ws.Range("C260").Value = "[name=`"Mousse`"]'With the guidance of the divine light, both will have bright futures.' `n"
